# "Refined metadata to be additional tab"
# 1. Re-stamp the "time_taken" query timestamps on the existing "data" sheet
#    (same rows/columns, new values captured at the later query run).
# 2. Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (name/id/version/version-created/query-time/request URL).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

$dataSheet.Range("F2").Value  = "2021-10-05 14:35:11.837741"
$dataSheet.Range("F3").Value  = "2021-10-05 14:35:11.837749"
$dataSheet.Range("F4").Value  = "2021-10-05 14:35:11.837752"
$dataSheet.Range("F5").Value  = "2021-10-05 14:35:11.837754"
$dataSheet.Range("F6").Value  = "2021-10-05 14:35:11.837757"
$dataSheet.Range("F7").Value  = "2021-10-05 14:35:11.837760"
$dataSheet.Range("F8").Value  = "2021-10-05 14:35:11.837763"
$dataSheet.Range("F9").Value  = "2021-10-05 14:35:11.837765"
$dataSheet.Range("F10").Value = "2021-10-05 14:35:11.837768"
$dataSheet.Range("F11").Value = "2021-10-05 14:35:11.837770"
$dataSheet.Range("F12").Value = "2021-10-05 14:35:11.837773"
$dataSheet.Range("F13").Value = "2021-10-05 14:35:11.837775"
$dataSheet.Range("F14").Value = "2021-10-05 14:35:11.837778"
$dataSheet.Range("F15").Value = "2021-10-05 14:35:11.837780"
$dataSheet.Range("F16").Value = "2021-10-05 14:35:11.837782"
$dataSheet.Range("F17").Value = "2021-10-05 14:35:11.837785"
$dataSheet.Range("F18").Value = "2021-10-05 14:35:11.837788"
$dataSheet.Range("F19").Value = "2021-10-05 14:35:11.837790"
$dataSheet.Range("F20").Value = "2021-10-05 14:35:11.837793"
$dataSheet.Range("F21").Value = "2021-10-05 14:35:11.837795"
$dataSheet.Range("F22").Value = "2021-10-05 14:35:11.837798"
$dataSheet.Range("F23").Value = "2021-10-05 14:35:11.837801"
$dataSheet.Range("F24").Value = "2021-10-05 14:35:11.837803"
$dataSheet.Range("F25").Value = "2021-10-05 14:35:11.837806"
$dataSheet.Range("F26").Value = "2021-10-05 14:35:11.837809"
$dataSheet.Range("F27").Value = "2021-10-05 14:35:11.837811"
$dataSheet.Range("F28").Value = "2021-10-05 14:35:11.837814"
$dataSheet.Range("F29").Value = "2021-10-05 14:35:11.837816"
$dataSheet.Range("F30").Value = "2021-10-05 14:35:11.837819"
$dataSheet.Range("F31").Value = "2021-10-05 14:35:11.837821"
$dataSheet.Range("F32").Value = "2021-10-05 14:35:11.837824"
$dataSheet.Range("F33").Value = "2021-10-05 14:35:11.837826"
$dataSheet.Range("F34").Value = "2021-10-05 14:35:11.837829"
$dataSheet.Range("F35").Value = "2021-10-05 14:35:11.837831"
$dataSheet.Range("F36").Value = "2021-10-05 14:35:11.837834"
$dataSheet.Range("F37").Value = "2021-10-05 14:35:11.837836"
$dataSheet.Range("F38").Value = "2021-10-05 14:35:11.837839"
$dataSheet.Range("F39").Value = "2021-10-05 14:35:11.837841"
$dataSheet.Range("F40").Value = "2021-10-05 14:35:11.837844"
$dataSheet.Range("F41").Value = "2021-10-05 14:35:11.837846"
$dataSheet.Range("F42").Value = "2021-10-05 14:35:11.837849"

# New "metadata" tab, placed right after "data" (sheetId 2, r:id rId2).
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Reuse the bold/bordered header style already defined in the workbook
# (style index used by the "data" header row) instead of authoring a new one.
$dataSheet.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)
$newSheet.Range("A2").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Paroxysmal Dyskinesia"
$newSheet.Range("C2").Value = 259
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.95"
$newSheet.Range("E2").Value = "2021-04-28T23:28:31.541048Z"
$newSheet.Range("F2").Value = "2021-10-05 14:35:11.834106"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/259/?format=json"

# Keep "data" as the active/visible sheet, as it was before this edit.
$dataSheet.Activate()
